$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new Y:AN columns for rows 1-16. Each cell references the
# corresponding "fraction" row 20 rows below (rows 21-36) and scales it up
# by 100000, e.g. Y1 = Y21*100000, Z1 = Z21*100000, ... AN16 = AN36*100000.
$ws.Range("Y1:AN1").Formula = "=Y21*100000"
$ws.Range("Y2:AN16").Formula = "=Y22*100000"

# Restore the selection/active-cell state recorded for the sheet.
$ws.Range("Y1:AN16").Select()
